$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(3).Delete()
$ws.Range("A26:A28").EntireRow.Delete()
